# Update "想去人数" (F column) figures across sheets, matching the
# upstream data refresh captured in the commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1894
$ws1.Range("F3").Value  = 1520
$ws1.Range("F4").Value  = 887
$ws1.Range("F5").Value  = 783
$ws1.Range("F6").Value  = 13372
$ws1.Range("F7").Value  = 13239
$ws1.Range("F11").Value = 566
$ws1.Range("F13").Value = 684
$ws1.Range("F14").Value = 2101
$ws1.Range("F15").Value = 66
$ws1.Range("F16").Value = 46
$ws1.Range("F19").Value = 400
$ws1.Range("F20").Value = 268
$ws1.Range("F21").Value = 292
$ws1.Range("F23").Value = 763

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 20
$ws2.Range("F8").Value = 11

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 196

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 196
$ws4.Range("F3").Value  = 1894
$ws4.Range("F4").Value  = 1520
$ws4.Range("F5").Value  = 887
$ws4.Range("F7").Value  = 783
$ws4.Range("F8").Value  = 13373
$ws4.Range("F9").Value  = 13239
$ws4.Range("F13").Value = 566
$ws4.Range("F15").Value = 684
$ws4.Range("F17").Value = 20
$ws4.Range("F18").Value = 2101
$ws4.Range("F19").Value = 66
$ws4.Range("F20").Value = 46
$ws4.Range("F26").Value = 400
$ws4.Range("F27").Value = 268
$ws4.Range("F28").Value = 292
$ws4.Range("F30").Value = 763
$ws4.Range("F32").Value = 11
